$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.612.09'
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").Value = '1.557.16'
$ws.Range("E3").Value = '  -1.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.485'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.56'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.99%  '
$ws.Range("E9").Value = '  -0.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0894'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").Value = '1.780.72'
$ws.Range("E12").Value = '  -1.10%  '
$ws.Range("D13").Value = '1.554.23'
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("D14").Value = '28.627.87'
$ws.Range("E14").Value = '  +1.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.513'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = '0.0₃0672'
$ws.Range("E20").Value = '  -2.12%  '
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("E22").Value = '  -1.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.59'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.40%  '
$ws.Range("E27").Value = '  -0.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0459'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.52%  '
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("E32").Value = '  -1.00%  '
$ws.Range("D33").Value = '1.391.81'
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.99'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.39%  '
$ws.Range("E35").Value = '  -2.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.65'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("E38").Value = '  -3.47%  '
$ws.Range("E39").Value = '  -0.68%  '
$ws.Range("E40").Value = '  +2.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.517'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.775'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0463'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.30'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.54%  '
$ws.Range("D47").Value = '1.693.25'
$ws.Range("E47").Value = '  -1.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.868'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '43.57'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '85.25'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").Value = '0.0₆0102'
$ws.Range("E51").Value = '  -0.51%  '
